# Scheduled market-data refresh: update profit calculation columns (H,I,J,K,L,M,N)
# across all job sheets, per upstream price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2834.3333
$ws.Range("J40").Value = 2923.6
$ws.Range("L40").Value = 2923.6
$ws.Range("N40").Value = -3273.6
$ws.Range("H47").Value = 19666.334
$ws.Range("I47").Value = 21999.5
$ws.Range("K47").Value = 21999.5
$ws.Range("M47").Value = -21027.5
$ws.Range("H86").Value = 4060.875
$ws.Range("I86").Value = 2830
$ws.Range("J86").Value = 4799.4
$ws.Range("K86").Value = 2830
$ws.Range("L86").Value = 4799.4
$ws.Range("M86").Value = -1707
$ws.Range("N86").Value = -7045.4
$ws.Range("H89").Value = 4060.875
$ws.Range("I89").Value = 2830
$ws.Range("J89").Value = 4799.4
$ws.Range("K89").Value = 14150
$ws.Range("L89").Value = 23997
$ws.Range("M89").Value = -8534
$ws.Range("N89").Value = -35229
$ws.Range("H92").Value = 1618.1
$ws.Range("I92").Value = 724.73334
$ws.Range("K92").Value = 724.73334
$ws.Range("M92").Value = 523.26666
$ws.Range("H99").Value = 2337.8333
$ws.Range("I99").Value = 342.66666
$ws.Range("J99").Value = 4333
$ws.Range("K99").Value = 1027.99998
$ws.Range("L99").Value = 12999
$ws.Range("M99").Value = 470.0000199999999
$ws.Range("N99").Value = -15995
$ws.Range("H101").Value = 3856.2856
$ws.Range("J101").Value = 7879.6
$ws.Range("L101").Value = 23638.8
$ws.Range("N101").Value = -26882.8
$ws.Range("H116").Value = 13701.807
$ws.Range("J116").Value = 11896.28
$ws.Range("L116").Value = 11896.28
$ws.Range("N116").Value = -18780.28
$ws.Range("H138").Value = 4041.158
$ws.Range("J138").Value = 3077.111
$ws.Range("L138").Value = 9231.332999999999
$ws.Range("N138").Value = -19511.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 64999.668
$ws.Range("J44").Value = 64999.668
$ws.Range("L44").Value = 64999.668
$ws.Range("N44").Value = -65975.66800000001
$ws.Range("H55").Value = 49998.5
$ws.Range("J55").Value = 49998.5
$ws.Range("L55").Value = 49998.5
$ws.Range("N55").Value = -50628.5
$ws.Range("H102").Value = 1709.6154
$ws.Range("I102").Value = 1709.6154
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1709.6154
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -87.61539999999991
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1655.8
$ws.Range("I132").Value = 1308.0526
$ws.Range("J132").Value = 2757
$ws.Range("K132").Value = 3924.1578
$ws.Range("L132").Value = 8271
$ws.Range("M132").Value = -1394.1578
$ws.Range("N132").Value = -13331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4728.5557
$ws.Range("I31").Value = 3019.6667
$ws.Range("J31").Value = 5583
$ws.Range("K31").Value = 3019.6667
$ws.Range("L31").Value = 5583
$ws.Range("M31").Value = -2724.6667
$ws.Range("N31").Value = -6173
$ws.Range("H34").Value = 4728.5557
$ws.Range("I34").Value = 3019.6667
$ws.Range("J34").Value = 5583
$ws.Range("K34").Value = 3019.6667
$ws.Range("L34").Value = 5583
$ws.Range("M34").Value = -2817.6667
$ws.Range("N34").Value = -5987
$ws.Range("H99").Value = 2942.5833
$ws.Range("I99").Value = 2626.7778
$ws.Range("J99").Value = 3890
$ws.Range("K99").Value = 2626.7778
$ws.Range("L99").Value = 3890
$ws.Range("M99").Value = -1128.7778
$ws.Range("N99").Value = -6886
$ws.Range("H107").Value = 4528.533
$ws.Range("I107").Value = 5557.364
$ws.Range("J107").Value = 1699.25
$ws.Range("K107").Value = 5557.364
$ws.Range("L107").Value = 1699.25
$ws.Range("M107").Value = -3637.364
$ws.Range("N107").Value = -5539.25
$ws.Range("H126").Value = 2942.5833
$ws.Range("I126").Value = 2626.7778
$ws.Range("J126").Value = 3890
$ws.Range("K126").Value = 7880.3334
$ws.Range("L126").Value = 11670
$ws.Range("M126").Value = -5410.3334
$ws.Range("N126").Value = -16610
$ws.Range("H134").Value = 4847.029
$ws.Range("I134").Value = 4391.129
$ws.Range("J134").Value = 8380.25
$ws.Range("K134").Value = 13173.387
$ws.Range("L134").Value = 25140.75
$ws.Range("M134").Value = -10638.387
$ws.Range("N134").Value = -30210.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 335.5
$ws.Range("I18").Value = 150
$ws.Range("J18").Value = 428.25
$ws.Range("K18").Value = 450
$ws.Range("L18").Value = 1284.75
$ws.Range("M18").Value = -281
$ws.Range("N18").Value = -1622.75
$ws.Range("H136").Value = 18506.059
$ws.Range("I136").Value = 5594.125
$ws.Range("J136").Value = 29983.334
$ws.Range("K136").Value = 16782.375
$ws.Range("L136").Value = 89950.00199999999
$ws.Range("M136").Value = -11682.375
$ws.Range("N136").Value = -100150.002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 56662
$ws.Range("J15").Value = 56662
$ws.Range("L15").Value = 56662
$ws.Range("N15").Value = -57238
$ws.Range("H81").Value = 56662
$ws.Range("J81").Value = 56662
$ws.Range("L81").Value = 56662
$ws.Range("N81").Value = -58658
$ws.Range("H84").Value = 56662
$ws.Range("J84").Value = 56662
$ws.Range("L84").Value = 169986
$ws.Range("N84").Value = -179970
$ws.Range("H107").Value = 330.4
$ws.Range("I107").Value = 341.375
$ws.Range("J107").Value = 286.5
$ws.Range("K107").Value = 341.375
$ws.Range("L107").Value = 286.5
$ws.Range("M107").Value = 1578.625
$ws.Range("N107").Value = -4126.5
$ws.Range("H122").Value = 3681.2
$ws.Range("I122").Value = 3302.3333
$ws.Range("J122").Value = 4249.5
$ws.Range("K122").Value = 9906.999899999999
$ws.Range("L122").Value = 12748.5
$ws.Range("M122").Value = -7456.999899999999
$ws.Range("N122").Value = -17648.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2037.1923
$ws.Range("I68").Value = 2070.5264
$ws.Range("K68").Value = 2070.5264
$ws.Range("M68").Value = -1321.5264
$ws.Range("H71").Value = 2037.1923
$ws.Range("I71").Value = 2070.5264
$ws.Range("K71").Value = 10352.632
$ws.Range("M71").Value = -6608.632000000001
$ws.Range("H80").Value = 29166.666
$ws.Range("J80").Value = 29166.666
$ws.Range("L80").Value = 29166.666
$ws.Range("N80").Value = -31412.666
$ws.Range("H82").Value = 3827.7222
$ws.Range("I82").Value = 1763.6364
$ws.Range("J82").Value = 7071.2856
$ws.Range("K82").Value = 1763.6364
$ws.Range("L82").Value = 7071.2856
$ws.Range("M82").Value = -1402.6364
$ws.Range("N82").Value = -7793.2856
$ws.Range("H83").Value = 29166.666
$ws.Range("J83").Value = 29166.666
$ws.Range("L83").Value = 87499.99800000001
$ws.Range("N83").Value = -98731.99800000001
$ws.Range("H85").Value = 3827.7222
$ws.Range("I85").Value = 1763.6364
$ws.Range("J85").Value = 7071.2856
$ws.Range("K85").Value = 1763.6364
$ws.Range("L85").Value = 7071.2856
$ws.Range("M85").Value = -515.6364000000001
$ws.Range("N85").Value = -9567.285599999999
$ws.Range("H132").Value = 2883.303
$ws.Range("I132").Value = 1918.9584
$ws.Range("J132").Value = 5454.8887
$ws.Range("K132").Value = 5756.8752
$ws.Range("L132").Value = 16364.6661
$ws.Range("M132").Value = -3226.8752
$ws.Range("N132").Value = -21424.6661
$ws.Range("H136").Value = 3007.4688
$ws.Range("I136").Value = 1200.5
$ws.Range("J136").Value = 3609.7917
$ws.Range("K136").Value = 3601.5
$ws.Range("L136").Value = 10829.3751
$ws.Range("M136").Value = -1051.5
$ws.Range("N136").Value = -15929.3751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 54069.125
$ws.Range("J45").Value = 54069.125
$ws.Range("L45").Value = 54069.125
$ws.Range("N45").Value = -55051.125
$ws.Range("H62").Value = 7517
$ws.Range("I62").Value = 8830.799999999999
$ws.Range("J62").Value = 5874.75
$ws.Range("K62").Value = 8830.799999999999
$ws.Range("L62").Value = 5874.75
$ws.Range("M62").Value = -8206.799999999999
$ws.Range("N62").Value = -7122.75
$ws.Range("H65").Value = 7517
$ws.Range("I65").Value = 8830.799999999999
$ws.Range("J65").Value = 5874.75
$ws.Range("K65").Value = 44154
$ws.Range("L65").Value = 29373.75
$ws.Range("M65").Value = -41034
$ws.Range("N65").Value = -35613.75
$ws.Range("H107").Value = 1551.3572
$ws.Range("I107").Value = 988.7273
$ws.Range("J107").Value = 3614.3333
$ws.Range("K107").Value = 2966.1819
$ws.Range("L107").Value = 10842.9999
$ws.Range("M107").Value = -1046.1819
$ws.Range("N107").Value = -14682.9999
$ws.Range("H126").Value = 1926
$ws.Range("I126").Value = 889
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 2667
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -197
$ws.Range("N126").Value = -16940

